# Update the task name for the Scoreboard Code Jam row and refresh the
# active selection, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix task title: "CodeJam "Scoreboard"" -> "Code Jam "Scoreboard""
$ws.Range("A8").Value = 'Code Jam "Scoreboard"'

# Move the active selection to A8
$ws.Range("A8").Select()
